$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9 and 10 swap coin identity (Cardano <-> Solana), with new price/volume
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"

# Updated Price (D) and Volume(1h) (E) values for each row.
# Price column is forced to Text (NumberFormat "@") before assignment so
# numeric-looking strings (e.g. "0.6750", "1.460") keep their exact textual
# representation instead of being parsed into numbers and losing trailing zeros.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.346.11"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.840.32"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.37"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6283"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07441"
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.99"
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2891"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07729"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.849.15"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.963"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6750"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001026"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.56"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.231"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.336.57"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.35"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.335"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.83"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.464"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1346"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07224"
$ws.Range("E28").Value = "  +12.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.460"
$ws.Range("E29").Value = "  +5.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.479"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.038"
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.033"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.821"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.138"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6946"
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.807"
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.835"
$ws.Range("E39").Value = "  +3.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.233.40"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9275"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9995"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.991.74"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.44"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.34"
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.707"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.947"
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.907"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1138"
$ws.Range("E50").Value = "  -3.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3902"
$ws.Range("E51").Value = "  -0.97%  "
